$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('E2').Value = '2026-02-10 21:18:35'
$ws.Range('E3').Value = '2026-02-10 21:18:38'
$ws.Range('E4').Value = '2026-02-10 21:18:40'
$ws.Range('E5').Value = '2026-02-10 21:18:43'
$ws.Range('E6').Value = '2026-02-10 21:18:45'
$ws.Range('I6').Value = '1.3 mm'
$ws.Range('J6').Value = '1004.1 hPa'
$ws.Range('O6').Value = '9.8 °C'
$ws.Range('E7').Value = '2026-02-10 21:18:48'
$ws.Range('J7').Value = '1004.4 hPa'
$ws.Range('E8').Value = '2026-02-10 21:18:50'
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = '83%'
$ws.Range('E9').Value = '2026-02-10 21:18:53'
$ws.Range('I9').Value = '3.3 mm'
$ws.Range('L9').Value = '9.0 km/h - 146º 20:41 TU'
$ws.Range('E10').Value = '2026-02-10 21:18:55'
$ws.Range('I10').Value = '1.1 mm'
$ws.Range('E11').Value = '2026-02-10 21:18:58'
$ws.Range('E12').Value = '2026-02-10 21:19:00'
$ws.Range('I12').Value = '4.1 mm'
$ws.Range('E13').Value = '2026-02-10 21:19:02'
$ws.Range('I13').Value = '10.7 mm'
$ws.Range('E14').Value = '2026-02-10 21:19:05'
$ws.Range('O14').Value = '13.2 °C'
$ws.Range('E15').Value = '2026-02-10 21:19:08'
$ws.Range('I15').Value = '3.7 mm'
$ws.Range('E16').Value = '2026-02-10 21:19:10'
$ws.Range('I16').Value = '26.3 mm'
$ws.Range('E17').Value = '2026-02-10 21:19:13'
$ws.Range('E18').Value = '2026-02-10 21:19:16'
$ws.Range('L18').Value = '18.4 km/h - 28º 20:53 TU'
$ws.Range('O18').Value = '10.4 °C'
$ws.Range('E19').Value = '2026-02-10 21:19:18'
$ws.Range('I19').Value = '0.5 mm'
$ws.Range('E20').Value = '2026-02-10 21:19:21'
$ws.Range('I20').Value = '11.4 mm'
$ws.Range('E21').Value = '2026-02-10 21:19:24'
$ws.Range('J21').Value = '1006.0 hPa'
$ws.Range('E22').Value = '2026-02-10 21:19:26'
$ws.Range('E23').Value = '2026-02-10 21:19:29'
$ws.Range('I23').Value = '27.1 mm'
$ws.Range('E24').Value = '2026-02-10 21:19:32'
$ws.Range('O24').Value = '11.2 °C'
$ws.Range('E25').Value = '2026-02-10 21:19:34'
$ws.Range('I25').Value = '22.3 mm'
$ws.Range('E26').Value = '2026-02-10 21:19:37'
$ws.Range('I26').Value = '0.7 mm'
$ws.Range('L26').Value = '22.0 km/h - 347º 20:37 TU'
$ws.Range('E27').Value = '2026-02-10 21:19:39'
$ws.Range('I27').Value = '12.3 mm'
$ws.Range('E28').Value = '2026-02-10 21:19:42'
$ws.Range('H28').NumberFormat = '@'
$ws.Range('H28').Value = '83%'
$ws.Range('J28').Value = '1004.1 hPa'
$ws.Range('E29').Value = '2026-02-10 21:19:45'
$ws.Range('L29').Value = '21.2 km/h - 349º 20:45 TU'
$ws.Range('E30').Value = '2026-02-10 21:19:47'
$ws.Range('I30').Value = '1.5 mm'
$ws.Range('L30').Value = '30.2 km/h - 353º 20:42 TU'
$ws.Range('E31').Value = '2026-02-10 21:19:50'
$ws.Range('H31').NumberFormat = '@'
$ws.Range('H31').Value = '82%'
$ws.Range('I31').Value = '3.3 mm'
$ws.Range('J31').Value = '1003.4 hPa'
$ws.Range('O31').Value = '10.5 °C'
$ws.Range('E32').Value = '2026-02-10 21:19:53'
$ws.Range('H32').NumberFormat = '@'
$ws.Range('H32').Value = '89%'
$ws.Range('E33').Value = '2026-02-10 21:19:55'
$ws.Range('I33').Value = '12.2 mm'
$ws.Range('L33').Value = '14.8 km/h - 330º 20:59 TU'
$ws.Range('E34').Value = '2026-02-10 21:19:57'
$ws.Range('I34').Value = '15.9 mm'
$ws.Range('E35').Value = '2026-02-10 21:20:00'
$ws.Range('J35').Value = '1004.9 hPa'
$ws.Range('O35').Value = '13.1 °C'
$ws.Range('E36').Value = '2026-02-10 21:20:03'
$ws.Range('I36').Value = '5.3 mm'
$ws.Range('J36').Value = '1004.3 hPa'
$ws.Range('E37').Value = '2026-02-10 21:20:06'
$ws.Range('I37').Value = '3.6 mm'
$ws.Range('O37').Value = '6.8 °C'
$ws.Range('E38').Value = '2026-02-10 21:20:08'
$ws.Range('O38').Value = '11.1 °C'
$ws.Range('E39').Value = '2026-02-10 21:20:11'
$ws.Range('I39').Value = '11.7 mm'
$ws.Range('E40').Value = '2026-02-10 21:20:13'
$ws.Range('H40').NumberFormat = '@'
$ws.Range('H40').Value = '91%'
$ws.Range('I40').Value = '14.8 mm'
$ws.Range('E41').Value = '2026-02-10 21:20:16'
$ws.Range('J41').Value = '1004.5 hPa'
$ws.Range('K41').Value = '9.5 MJ/m2'
$ws.Range('L41').Value = '40.0 km/h - 274º 20:49 TU'
$ws.Range('M41').Value = '21.6 °C 20:52 TU'
$ws.Range('O41').Value = '14.6 °C'
$ws.Range('E42').Value = '2026-02-10 21:20:19'
$ws.Range('I42').Value = '0.9 mm'
$ws.Range('O42').Value = '10.4 °C'
$ws.Range('E43').Value = '2026-02-10 21:20:21'
$ws.Range('O43').Value = '9.9 °C'
$ws.Range('E44').Value = '2026-02-10 21:20:23'
$ws.Range('I44').Value = '29.4 mm'
$ws.Range('E45').Value = '2026-02-10 21:20:26'
$ws.Range('E46').Value = '2026-02-10 21:20:29'
$ws.Range('O46').Value = '14.8 °C'
